$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new Time Log entry (row 25): 2/22/2019, 5.5 hrs, task note ---
$ws.Range("A25").Value = 43518
$ws.Range("A25").NumberFormat = "d-mmm"
$ws.Range("B25").Value = 5.5
$ws.Range("D25").Value = "Indie Project: MVP ERD, populated lookup tables, tried to figure mysqldump (saved creation & insert files separately for now), created User entity"
$ws.Rows.Item(25).RowHeight = 30

# --- Append a follow-up note below the existing bottom note (row 32 stays, new row 33) ---
$ws.Range("D33").Value = "Now I have a path to get at mysqldump but I have an access problem for writing the dump to the locations I choose"

# --- Update view state: scroll so row 16 is at top, select D34 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("D34").Select() | Out-Null

# --- Maximize the application window (reflects the saved workbook view geometry) ---
$excel.WindowState = -4137

Write-Host "Edit complete"
